$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Text replacements (run formatting is preserved automatically by Find)
# ---------------------------------------------------------------------------

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Nanotechnology in Medicine: Shaping the Future of Healthcare" "The Symphony of Atoms: Exploring the Marvels of Chemistry"

Replace-Text " Alex Morgan" " Emily Carter"

# "alex" (standalone run) must not match the "Alex" inside " Alex Morgan" above,
# which has already been replaced by this point, so a case-sensitive / whole
# word match is safe here.
$d.Content.Find.Execute("alex", $true, $true, $false, $false, $false, $true, 1, $false, "Emily", 2) | Out-Null

Replace-Text "morgan@healthsciences" "Carter@edumail"

Replace-Text "Nanotechnology, the manipulation of matter at a nanoscale (1-100 nanometers), holds immense potential to revolutionize various fields, including medicine" "The universe we inhabit is a symphony of atoms, an intricate dance of particles that orchestrate the very essence of life itself"

Replace-Text " By harnessing the unique properties of materials at this scale, scientists can engineer novel medical devices, treatments, and drug delivery systems with unprecedented precision and efficacy" " Chemistry, the study of the properties, behavior, and interactions of matter, unveils the secrets of this atomic ballet, guiding us into the depths of the natural world"

Replace-Text " In this essay, we will delve into the transformative applications of nanotechnology in medicine and explore how it is poised to reshape the future of healthcare" " This captivating field unravels the mechanisms behind everything from the mundane to the miraculous, the rusting of iron to the blooming of a flower, the marvels of medicine to the complexities of life's processes"

Replace-Text "Nanotechnology offers the potential to overcome the limitations of conventional medical approaches by enabling targeted, non-invasive interventions at the cellular and molecular level" "Chemistry, with its myriad elements and compounds, serves as a potent lens through which we can explore the intricate tapestry of nature"

Replace-Text " For instance, nanoscale drug delivery systems can be designed to specifically target diseased cells while sparing healthy tissues, minimizing side effects and improving treatment outcomes" " It enables us to unravel the molecular underpinnings of matter, understanding the structure and composition of substances and deciphering the intricate interactions that govern their behavior"

Replace-Text " Additionally, the development of nano-enabled sensors and imaging techniques allows for real-time monitoring of physiological processes, enabling early detection and intervention in diseases" " From the vastness of the cosmos to the microscopic realm of atoms and molecules, chemistry offers a pathway to comprehending the complex realities that surround us"

Replace-Text "Furthermore, nanotechnology can empower the development of innovative, personalized treatments tailored to individual genetic profiles and medical conditions" "The study of chemistry empowers us to not only comprehend the world but also to harness its powers for the betterment of humankind"

Replace-Text " By harnessing the ability to manipulate matter at the molecular level, researchers can engineer therapies that precisely target specific disease pathways or genetic mutations" " Through chemical advancements, we have witnessed the development of life-saving medicines, innovative materials, and sustainable energy solutions"

Replace-Text " This approach holds the promise of transformative outcomes for patients with complex or currently incurable diseases" " Chemistry continues to unveil profound implications for tackling societal challenges, from mitigating climate change to ensuring food security, demonstrating its profound impact on shaping our future"

Replace-Text "Nanotechnology in medicine is a rapidly evolving field with the potential to revolutionize healthcare" "Chemistry is the exploration of the composition, behavior, and interactions of matter"

Replace-Text " It provides a powerful platform for developing targeted drug delivery systems, enhancing diagnostic capabilities, and engineering personalized treatments" " It unravels the atomic ballet that orchestrates the world around us, connecting phenomena from the mundane to the miraculous"

Replace-Text " By exploiting the unique properties of materials at the nanoscale, researchers can design innovative solutions that can precisely target diseased cells, minimize side effects, enable early disease detection, and empower personalized medicine" " This field unveils the intricacies of matter, structure, and reactivity, offering a pathway to understanding the natural world"

Replace-Text " As nanotechnology continues to advance, we can anticipate groundbreaking advancements in healthcare, leading to improved patient outcomes and a healthier future" " Chemistry empowers us to comprehend and manipulate the material realm, leading to advancements in medicine, materials science, and energy solutions, ultimately shaping our future and addressing societal challenges"

# ---------------------------------------------------------------------------
# 2. Insert the two brand-new runs (sentence + its period) right after the
#    "...life's processes" sentence's trailing period, before the following
#    line break.
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("life's processes", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
# the period run immediately follows the sentence we just found
$rng.MoveEnd(1, 1) | Out-Null
$rng.Collapse(0)

$rng.InsertAfter(" Chemistry weaves together the fabric of our material existence, forging connections between diverse phenomena and illuminating the fundamental principles governing our world")
$rng.Font.Name = "Times New Roman"
$rng.Font.Color = 0
$rng.Font.Size = 12

$rng.Collapse(0)
$rng.InsertAfter(".")
$rng.Font.Name = "Times New Roman"
$rng.Font.Color = 0
$rng.Font.Size = 12

# ---------------------------------------------------------------------------
# 3. Append a new, empty paragraph at the very end of the document.
# ---------------------------------------------------------------------------

$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 4. Font rename pass: TimesNewToman -> Times New Roman, across the whole
#    document (the two brand-new runs above are already correct).
# ---------------------------------------------------------------------------

$full = $d.Range(0, $d.Content.End)
$full.Font.Name = "Times New Roman"
